$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 44224, "Primera", 200, 16500, 17000, 16750, 931),
    @(3, 44224, "Segunda", 200, 13500, 14000, 13750, 764),
    @(4, 44385, "Especial", 120, 14000, 14500, 14250, 792),
    @(5, 44385, "Primera", 300, 11000, 11500, 11250, 625),
    @(6, 44385, "Segunda", 240, 8000, 8500, 8250, 458),
    @(7, 44385, "Tercera", 120, 5000, 5500, 5250, 292),
    @(8, 44273, "Especial", 160, 12500, 13000, 12750, 708),
    @(9, 44273, "Primera", 240, 10500, 11000, 10750, 597),
    @(10, 44273, "Segunda", 200, 8500, 9000, 8750, 486),
    @(11, 44308, "Especial", 300, 15500, 16000, 15750, 875),
    @(12, 44308, "Primera", 240, 13500, 14000, 13750, 764),
    @(13, 44308, "Segunda", 200, 10500, 11000, 10750, 597),
    @(14, 44335, "Especial", 240, 19500, 20000, 19750, 1097),
    @(15, 44335, "Primera", 200, 17500, 18000, 17750, 986),
    @(16, 44335, "Segunda", 160, 12500, 13000, 12750, 708),
    @(17, 44272, "Especial", 160, 12500, 13000, 12750, 708),
    @(18, 44272, "Primera", 300, 10500, 11000, 10750, 597),
    @(19, 44272, "Segunda", 240, 8500, 9000, 8750, 486),
    @(20, 44280, "Especial", 240, 12500, 13000, 12750, 708),
    @(21, 44280, "Primera", 240, 10500, 11000, 10750, 597),
    @(22, 44280, "Segunda", 300, 8500, 9000, 8750, 486),
    @(23, 44286, "Especial", 700, 12500, 13000, 12750, 708),
    @(24, 44286, "Primera", 500, 10500, 11000, 10750, 597),
    @(25, 44286, "Segunda", 300, 8500, 9000, 8750, 486),
    @(26, 44294, "Especial", 360, 12500, 13000, 12750, 708),
    @(27, 44294, "Primera", 240, 10500, 11000, 10750, 597),
    @(28, 44294, "Segunda", 240, 8500, 9000, 8750, 486),
    @(29, 44293, "Especial", 400, 12500, 13000, 12750, 708),
    @(30, 44293, "Primera", 508, 10500, 11000, 10746, 597),
    @(31, 44293, "Segunda", 400, 8500, 9000, 8750, 486),
    @(32, 44279, "Especial", 200, 12500, 13000, 12750, 708),
    @(33, 44279, "Primera", 240, 10500, 11000, 10750, 597),
    @(34, 44279, "Segunda", 240, 8500, 9000, 8750, 486),
    @(35, 44384, "Especial", 160, 14500, 15000, 14750, 819),
    @(36, 44384, "Primera", 500, 11500, 12000, 11750, 653),
    @(37, 44384, "Segunda", 400, 8500, 9000, 8750, 486),
    @(38, 44384, "Tercera", 300, 5500, 6000, 5750, 319)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 4).Value = $row[1]   # D: Fecha
    $ws.Cells.Item($r, 12).Value = $row[2]  # L: Calidad
    $ws.Cells.Item($r, 13).Value = $row[3]  # M: Volumen
    $ws.Cells.Item($r, 14).Value = $row[4]  # N: Precio minimo
    $ws.Cells.Item($r, 15).Value = $row[5]  # O: Precio maximo
    $ws.Cells.Item($r, 16).Value = $row[6]  # P: Precio promedio ponderado
    $ws.Cells.Item($r, 19).Value = $row[7]  # S: Precio $/Kg
}

$wb.Save()
